# Fix current and prior year headers in remaining QC town close workbooks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns H:M previously held "2023 ..." labels -> now "Prior Year ..."
$ws.Range("H1").Value = "Prior Year LMV"
$ws.Range("I1").Value = "Prior Year BMV"
$ws.Range("J1").Value = "Prior Year Total MV"
$ws.Range("K1").Value = "Prior Year LAV"
$ws.Range("L1").Value = "Prior Year BAV"
$ws.Range("M1").Value = "Prior Year Total AV"

# Columns N:S previously held "2024 ..." labels -> now "Curr. Year ..."
$ws.Range("N1").Value = "Curr. Year LMV"
$ws.Range("O1").Value = "Curr. Year BMV"
$ws.Range("P1").Value = "Curr. Year Total MV"
$ws.Range("Q1").Value = "Curr. Year LAV"
$ws.Range("R1").Value = "Curr. Year BAV"
$ws.Range("S1").Value = "Curr. Year Total AV"
